# Re-key the IFRS figures onto the restated (smaller) reporting unit
# and drop the rows/columns the restated sheet no longer carries.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 6045
$ws.Range("E2").Value = 684
$ws.Range("F2").Value = 684
$ws.Range("G2").Value = 734
$ws.Range("H2").Value = 575
$ws.Range("I2").Value = 575
$ws.Range("K2").Value = 6467
$ws.Range("L2").Value = 1477
$ws.Range("M2").Value = 4990
$ws.Range("N2").Value = 4990
$ws.Range("P2").Value = 52
$ws.Range("Q2").Value = 768
$ws.Range("R2").Value = -808
$ws.Range("S2").Value = -35
$ws.Range("T2").Value = 834
$ws.Range("U2").Value = -66
$ws.Range("V2").Value = 92
$ws.Range("W2").Value = 11.31
$ws.Range("X2").Value = 9.52
$ws.Range("Y2").Value = 12.28
$ws.Range("Z2").Value = 9.44
$ws.Range("AA2").Value = 29.61
$ws.Range("AB2").Value = 9130.120000000001
$ws.Range("AC2").Value = 5525
$ws.Range("AD2").Value = 12.22
$ws.Range("AE2").Value = 47909
$ws.Range("AF2").Value = 1.41
$ws.Range("AG2").Value = 450
$ws.Range("AH2").Value = 0.67
$ws.Range("AI2").Value = 8.140000000000001
$ws.Range("AJ2").Value = 10415000
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3
$ws.Range("D3").Value = 6613
$ws.Range("E3").Value = 827
$ws.Range("F3").Value = 827
$ws.Range("G3").Value = 880
$ws.Range("H3").Value = 689
$ws.Range("I3").Value = 689
$ws.Range("K3").Value = 7039
$ws.Range("L3").Value = 1442
$ws.Range("M3").Value = 5596
$ws.Range("N3").Value = 5596
$ws.Range("P3").Value = 52
$ws.Range("Q3").Value = 931
$ws.Range("R3").Value = -647
$ws.Range("S3").Value = -65
$ws.Range("T3").Value = 582
$ws.Range("U3").Value = 349
$ws.Range("V3").Value = 94
$ws.Range("W3").Value = 12.5
$ws.Range("X3").Value = 10.42
$ws.Range("Y3").Value = 13.01
$ws.Range("Z3").Value = 10.2
$ws.Range("AA3").Value = 25.77
$ws.Range("AB3").Value = 10335.81
$ws.Range("AC3").Value = 6614
$ws.Range("AD3").Value = 14.85
$ws.Range("AE3").Value = 53734
$ws.Range("AF3").Value = 1.83
$ws.Range("AG3").Value = 600
$ws.Range("AH3").Value = 0.61
$ws.Range("AI3").Value = 9.07
$ws.Range("AJ3").Value = 10415000
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 4
$ws.Range("D4").Value = 7142
$ws.Range("E4").Value = 904
$ws.Range("F4").Value = 904
$ws.Range("G4").Value = 886
$ws.Range("H4").Value = 672
$ws.Range("I4").Value = 672
$ws.Range("K4").Value = 8218
$ws.Range("L4").Value = 2074
$ws.Range("M4").Value = 6144
$ws.Range("N4").Value = 6144
$ws.Range("P4").Value = 52
$ws.Range("Q4").Value = 1489
$ws.Range("R4").Value = -1067
$ws.Range("S4").Value = -130
$ws.Range("T4").Value = 491
$ws.Range("U4").Value = 998
$ws.Range("V4").Value = 34
$ws.Range("W4").Value = 12.66
$ws.Range("X4").Value = 9.41
$ws.Range("Y4").Value = 11.45
$ws.Range("Z4").Value = 8.81
$ws.Range("AA4").Value = 33.75
$ws.Range("AB4").Value = 11476.6
$ws.Range("AC4").Value = 6452
$ws.Range("AD4").Value = 12.2
$ws.Range("AE4").Value = 58993
$ws.Range("AF4").Value = 1.33
$ws.Range("AG4").Value = 700
$ws.Range("AH4").Value = 0.89
$ws.Range("AI4").Value = 10.85
$ws.Range("AJ4").Value = 10415000
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# Row 5
$ws.Range("D5").Value = 7444
$ws.Range("E5").Value = 661
$ws.Range("F5").Value = 661
$ws.Range("G5").Value = 641
$ws.Range("H5").Value = 498
$ws.Range("I5").Value = 498
$ws.Range("K5").Value = 8752
$ws.Range("L5").Value = 1996
$ws.Range("M5").Value = 6756
$ws.Range("N5").Value = 6756
$ws.Range("P5").Value = 52
$ws.Range("Q5").Value = 767
$ws.Range("R5").Value = -754
$ws.Range("S5").Value = -121
$ws.Range("T5").Value = 779
$ws.Range("U5").Value = -12
$ws.Range("V5").Value = 28
$ws.Range("W5").Value = 8.880000000000001
$ws.Range("X5").Value = 6.69
$ws.Range("Y5").Value = 7.72
$ws.Range("Z5").Value = 5.87
$ws.Range("AA5").Value = 29.55
$ws.Range("AB5").Value = 12317.41
$ws.Range("AC5").Value = 4783
$ws.Range("AD5").Value = 13.99
$ws.Range("AE5").Value = 64864
$ws.Range("AF5").Value = 1.03
$ws.Range("AG5").Value = 700
$ws.Range("AH5").Value = 1.05
$ws.Range("AI5").Value = 14.63
$ws.Range("AJ5").Value = 10415000
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 7464
$ws.Range("E6").Value = 462
$ws.Range("F6").Value = 462
$ws.Range("G6").Value = 508
$ws.Range("H6").Value = 396
$ws.Range("I6").Value = 396
$ws.Range("K6").Value = 8605
$ws.Range("L6").Value = 1696
$ws.Range("M6").Value = 6909
$ws.Range("N6").Value = 6909
$ws.Range("P6").Value = 52
$ws.Range("Q6").Value = 530
$ws.Range("R6").Value = -652
$ws.Range("S6").Value = -72
$ws.Range("T6").Value = 588
$ws.Range("U6").Value = -58
$ws.Range("V6").Value = 31
$ws.Range("W6").Value = 6.19
$ws.Range("X6").Value = 5.31
$ws.Range("Y6").Value = 5.8
$ws.Range("Z6").Value = 4.56
$ws.Range("AA6").Value = 24.55
$ws.Range("AB6").Value = 13123.53
$ws.Range("AC6").Value = 3803
$ws.Range("AD6").Value = 10.68
$ws.Range("AE6").Value = 66337
$ws.Range("AF6").Value = 0.61
$ws.Range("AG6").Value = 700
$ws.Range("AH6").Value = 1.72
$ws.Range("AI6").Value = 18.41
$ws.Range("AJ6").Value = 10415000

# Forecast rows with no restated figures yet - keep only the labels
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
